# Update attendee counts on "展览" (Exhibitions) and "全部类型" (All types) sheets
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 2334
$wsExhibit.Range("F6").Value = 214
$wsExhibit.Range("F7").Value = 371

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 2334
$wsAll.Range("F6").Value = 214
$wsAll.Range("F9").Value = 371
